$wb = $excel.ActiveWorkbook

# --- Sheet 1: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.6179775280898876
$ws1.Range("C2").Value = 0.5668789808917197
$ws1.Range("D2").Value = 1
$ws1.Range("E2").Value = 0.7235772357723578
$ws1.Range("F2").Value = 0.8674463937621832
$ws1.Range("G2").Value = 0.9714525608732157
$ws1.Range("H2").Value = 0.8048331439633043
$ws1.Range("I2").Value = 534
$ws1.Range("J2").Value = 408
$ws1.Range("K2").Value = 126
$ws1.Range("L2").Value = 0

# --- Sheet 2: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")
$ws2.Range("B2").Value = 1
$ws2.Range("C2").Value = 0.2359550561797753
$ws2.Range("D2").Value = 0.3818181818181818

$ws2.Range("B3").Value = 0.5668789808917197
$ws2.Range("C3").Value = 1
$ws2.Range("D3").Value = 0.7235772357723578

$ws2.Range("B4").Value = 0.6179775280898876
$ws2.Range("C4").Value = 0.6179775280898876
$ws2.Range("D4").Value = 0.6179775280898876
$ws2.Range("E4").Value = 0.6179775280898876

$ws2.Range("B5").Value = 0.7834394904458599
$ws2.Range("C5").Value = 0.6179775280898876
$ws2.Range("D5").Value = 0.5526977087952698

$ws2.Range("B6").Value = 0.7834394904458598
$ws2.Range("C6").Value = 0.6179775280898876
$ws2.Range("D6").Value = 0.5526977087952698

# --- Sheet 3: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 126
$ws3.Range("C2").Value = 408

$ws3.Range("B3").Value = 0
$ws3.Range("C3").Value = 534
